# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) previously held a "Strike#" style count; this
# recomputes it to the real strikeout ("K") totals for each start and
# writes the corrected value back into column G for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 1
    4  = 5
    5  = 0
    6  = 3
    7  = 5
    8  = 5
    9  = 2
    10 = 6
    11 = 1
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    21 = 6
    22 = 0
    23 = 0
    24 = 1
    25 = 2
    26 = 2
    27 = 3
    28 = 3
    29 = 1
    30 = 2
    31 = 3
    32 = 0
    33 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
